# Generate Report for Handback
# Adds two new handback records to the report workbook:
#   1) 1e0cfe05-4bff-4dd8-b4d2-eda6a5e26823
#   2) e7b90767-91c5-4fc0-969b-599e9c24842e
# across all three sheets: Overview, zh-cn, de-de

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$hlColor = 15570276   # matches font color FF6495ED used by existing hyperlink cells
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---- record identifiers -------------------------------------------------
$md1 = "1e0cfe05-4bff-4dd8-b4d2-eda6a5e26823.md"
$md2 = "e7b90767-91c5-4fc0-969b-599e9c24842e.md"

$zh1xlf = "1e0cfe05-4bff-4dd8-b4d2-eda6a5e26823.ca3edbb921036ac4229229dc4bf882b64728fda9.zh-cn.xlf"
$zh2xlf = "e7b90767-91c5-4fc0-969b-599e9c24842e.6327b458e432c26df288783631f9233012e7d4c1.zh-cn.xlf"
$de1xlf = "1e0cfe05-4bff-4dd8-b4d2-eda6a5e26823.ca3edbb921036ac4229229dc4bf882b64728fda9.de-de.xlf"
$de2xlf = "e7b90767-91c5-4fc0-969b-599e9c24842e.6327b458e432c26df288783631f9233012e7d4c1.de-de.xlf"

$zh1HandoffDate = "2016-02-18 04:11:44"
$zh1HandbackDate = "2016-02-18 04:12:30"
$zh2HandoffDate = "2016-02-18 04:11:44"
$zh2HandbackDate = "2016-02-18 04:12:30"

$de1HandoffDate = "2016-02-18 04:11:57"
$de1HandbackDate = "2016-02-18 04:12:53"
$de2HandoffDate = "2016-02-18 04:11:57"
$de2HandbackDate = "2016-02-18 04:12:53"

$status = "Handed back: in sync with en-US"
$reason = "Include"

# =====================================================================
# Sheet "Overview": rows 6 and 7, columns A (File Name), B (zh-cn), C (de-de)
# =====================================================================
$ws1.Range("A6").Value = $md1
$ws1.Range("A7").Value = $md2

$ws1.Range("B6").Value = $status
$ws1.Range("C6").Value = $status
$ws1.Range("B7").Value = $status
$ws1.Range("C7").Value = $status

$ws1.Range("A6").Font.Color = $hlColor
$ws1.Range("A6").Font.Underline = 2
$ws1.Range("A7").Font.Color = $hlColor
$ws1.Range("A7").Font.Underline = 2

$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/$md1", "", "", $md1)
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/$md2", "", "", $md2)

# =====================================================================
# Sheet "zh-cn": rows 6 and 7
# Columns: A Source File Name, B Status, C Correspond Handoff File,
#          D Correspond Handoff Datetime, E Target File,
#          F Correspond Handback File, G Correspond Handback DateTime,
#          H Handoff Reason, I Dependency From
# =====================================================================
$ws2.Range("A6").Value = $md1
$ws2.Range("B6").Value = $status
$ws2.Range("C6").Value = $zh1xlf
$ws2.Range("D6").Value = $zh1HandoffDate
$ws2.Range("E6").Value = $md1
$ws2.Range("F6").Value = $zh1xlf
$ws2.Range("G6").Value = $zh1HandbackDate
$ws2.Range("H6").Value = $reason

$ws2.Range("A7").Value = $md2
$ws2.Range("B7").Value = $status
$ws2.Range("C7").Value = $zh2xlf
$ws2.Range("D7").Value = $zh2HandoffDate
$ws2.Range("E7").Value = $md2
$ws2.Range("F7").Value = $zh2xlf
$ws2.Range("G7").Value = $zh2HandbackDate
$ws2.Range("H7").Value = $reason

foreach ($ref in "A6","E6","A7","E7") {
    $ws2.Range($ref).Font.Color = $hlColor
    $ws2.Range($ref).Font.Underline = 2
}
foreach ($ref in "C6","F6","C7","F7") {
    $ws2.Range($ref).Font.Color = $hlColor
    $ws2.Range($ref).Font.Underline = 2
}
$ws2.Range("D6").NumberFormat = $dateFmt
$ws2.Range("D7").NumberFormat = $dateFmt

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ebc3f88674eb549ffd4e4fed464340339e8c6a0c/e2e/$md1", "", "", $md1)
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f643c153db888935ca2415b7c6605297bfa132b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zh1xlf", "", "", $zh1xlf)
$ws2.Hyperlinks.Add($ws2.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ebc3f88674eb549ffd4e4fed464340339e8c6a0c/e2e/$md1", "", "", $md1)
$ws2.Hyperlinks.Add($ws2.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e47ab174e7e5b89eb266d5d59a1c52b2c54bac41/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zh1xlf", "", "", $zh1xlf)

$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ebc3f88674eb549ffd4e4fed464340339e8c6a0c/e2e/$md2", "", "", $md2)
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f643c153db888935ca2415b7c6605297bfa132b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zh2xlf", "", "", $zh2xlf)
$ws2.Hyperlinks.Add($ws2.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ebc3f88674eb549ffd4e4fed464340339e8c6a0c/e2e/$md2", "", "", $md2)
$ws2.Hyperlinks.Add($ws2.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e47ab174e7e5b89eb266d5d59a1c52b2c54bac41/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zh2xlf", "", "", $zh2xlf)

# =====================================================================
# Sheet "de-de": rows 6 and 7 (same layout as zh-cn)
# =====================================================================
$ws3.Range("A6").Value = $md1
$ws3.Range("B6").Value = $status
$ws3.Range("C6").Value = $de1xlf
$ws3.Range("D6").Value = $de1HandoffDate
$ws3.Range("E6").Value = $md1
$ws3.Range("F6").Value = $de1xlf
$ws3.Range("G6").Value = $de1HandbackDate
$ws3.Range("H6").Value = $reason

$ws3.Range("A7").Value = $md2
$ws3.Range("B7").Value = $status
$ws3.Range("C7").Value = $de2xlf
$ws3.Range("D7").Value = $de2HandoffDate
$ws3.Range("E7").Value = $md2
$ws3.Range("F7").Value = $de2xlf
$ws3.Range("G7").Value = $de2HandbackDate
$ws3.Range("H7").Value = $reason

foreach ($ref in "A6","E6","A7","E7") {
    $ws3.Range($ref).Font.Color = $hlColor
    $ws3.Range($ref).Font.Underline = 2
}
foreach ($ref in "C6","F6","C7","F7") {
    $ws3.Range($ref).Font.Color = $hlColor
    $ws3.Range($ref).Font.Underline = 2
}
$ws3.Range("D6").NumberFormat = $dateFmt
$ws3.Range("D7").NumberFormat = $dateFmt

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d2590e99c3e7c94881c3fb4a014a506584edc2f/e2e/$md1", "", "", $md1)
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/895eb9b29c83e3db4a02ee7ce76aee4e446931fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$de1xlf", "", "", $de1xlf)
$ws3.Hyperlinks.Add($ws3.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d2590e99c3e7c94881c3fb4a014a506584edc2f/e2e/$md1", "", "", $md1)
$ws3.Hyperlinks.Add($ws3.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/80549469b1cc65a2487e522de43ffdc07430b42c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$de1xlf", "", "", $de1xlf)

$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d2590e99c3e7c94881c3fb4a014a506584edc2f/e2e/$md2", "", "", $md2)
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/895eb9b29c83e3db4a02ee7ce76aee4e446931fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$de2xlf", "", "", $de2xlf)
$ws3.Hyperlinks.Add($ws3.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d2590e99c3e7c94881c3fb4a014a506584edc2f/e2e/$md2", "", "", $md2)
$ws3.Hyperlinks.Add($ws3.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/80549469b1cc65a2487e522de43ffdc07430b42c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$de2xlf", "", "", $de2xlf)

Write-Host "Handback report rows appended for 1e0cfe05 and e7b90767 on all sheets."
